# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.086.96"
$ws.Range("E2").Value = "  +2.20%  "

# Row 3
$ws.Range("D3").Value = "2.276.18"
$ws.Range("E3").Value = "  +2.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "'318.89"
$ws.Range("E5").Value = "  +0.68%  "

# Row 6
$ws.Range("D6").Value = "'103.49"
$ws.Range("E6").Value = "  +4.83%  "

# Row 7
$ws.Range("D7").Value = "'0.586"
$ws.Range("E7").Value = "  +0.86%  "

# Row 8
$ws.Range("E8").Value = "  -0.21%  "

# Row 9
$ws.Range("D9").Value = "'0.573"
$ws.Range("E9").Value = "  +1.89%  "

# Row 10
$ws.Range("D10").Value = "'39.11"
$ws.Range("E10").Value = "  +6.47%  "

# Row 11
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  +1.34%  "

# Row 12
$ws.Range("E12").Value = "  +1.38%  "

# Row 13
$ws.Range("E13").Value = "  +1.98%  "

# Row 14
$ws.Range("D14").Value = "2.622.14"
$ws.Range("E14").Value = "  +2.33%  "

# Row 15
$ws.Range("D15").Value = "'0.877"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16
$ws.Range("D16").Value = "'14.57"
$ws.Range("E16").Value = "  +3.23%  "

# Row 17
$ws.Range("D17").Value = "2.262.92"
$ws.Range("E17").Value = "  +1.84%  "

# Row 18
$ws.Range("D18").Value = "43.977.65"
$ws.Range("E18").Value = "  +2.40%  "

# Row 19
$ws.Range("D19").Value = "'14.25"
$ws.Range("E19").Value = "  -3.64%  "

# Row 20
$ws.Range("D20").Value = "'0.0₂01000"
$ws.Range("E20").Value = "  +3.98%  "

# Row 21
$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  +3.21%  "

# Row 22
$ws.Range("D22").Value = "'66.24"
$ws.Range("E22").Value = "  +1.59%  "

# Row 23
$ws.Range("E23").Value = "  +0.56%  "

# Row 24
$ws.Range("D24").Value = "'237.69"
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("E25").Value = "  +2.90%  "

# Row 26
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("D27").Value = "'10.28"
$ws.Range("E27").Value = "  +1.70%  "

# Row 28
$ws.Range("D28").Value = "'38.85"
$ws.Range("E28").Value = "  +14.61%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").Value = "'6.57"
$ws.Range("E30").Value = "  +4.45%  "

# Row 31
$ws.Range("D31").Value = "'162.59"
$ws.Range("E31").Value = "  +4.77%  "

# Row 32
$ws.Range("D32").Value = "'20.49"
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("E33").Value = "  -1.08%  "

# Row 34
$ws.Range("D34").Value = "'2.71"
$ws.Range("E34").Value = "  -2.42%  "

# Row 35
$ws.Range("D35").Value = "'3.24"
$ws.Range("E35").Value = "  -0.57%  "

# Row 36
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  +0.48%  "

# Row 37
$ws.Range("E37").Value = "  -1.13%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.109"
$ws.Range("E38").Value = "  +5.08%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'4.52"
$ws.Range("E39").Value = "  +1.59%  "

# Row 40
$ws.Range("D40").Value = "'3.89"
$ws.Range("E40").Value = "  +3.76%  "

# Row 41
$ws.Range("D41").Value = "'15.57"
$ws.Range("E41").Value = "  +26.92%  "

# Row 42
$ws.Range("E42").Value = "  +0.44%  "

# Row 43
$ws.Range("E43").Value = "  -0.13%  "

# Row 44
$ws.Range("D44").Value = "1.770.59"
$ws.Range("E44").Value = "  -6.18%  "

# Row 45
$ws.Range("E45").Value = "  -0.08%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'85.14"
$ws.Range("E46").Value = "  -4.91%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'5.39"
$ws.Range("E47").Value = "  -1.31%  "

# Row 48
$ws.Range("D48").Value = "'8.88"
$ws.Range("E48").Value = "  +1.28%  "

# Row 49
$ws.Range("D49").Value = "'59.69"
$ws.Range("E49").Value = "  -1.61%  "

# Row 50
$ws.Range("D50").Value = "'74.83"
$ws.Range("E50").Value = "  -3.51%  "

# Row 51
$ws.Range("D51").Value = "'104.52"
$ws.Range("E51").Value = "  +3.45%  "

